$wb = $excel.ActiveWorkbook

# OFF sheet
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 252
$wsOff.Range("C2").Value = 158
$wsOff.Range("D2").Value = 49
$wsOff.Range("E2").Value = 25

# DEF sheet
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 231
$wsDef.Range("C2").Value = 173
$wsDef.Range("D2").Value = 58
$wsDef.Range("E2").Value = 36
